# The "language_metadata" table had two rows removed:
#   - row 10: A="en  (malformed/unterminated quote), B=0, C=1
#   - the old row 19: A=Yes, B=0, C=1
# and the count for the `"en"` row increased from 41 to 43.
#
# Deleting entire rows shifts everything below up, so we delete from the
# bottom first (old row 19) to keep the row 10 reference valid for the
# second delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("19").EntireRow.Delete() | Out-Null
$ws.Rows("10").EntireRow.Delete() | Out-Null

# After the deletes, the old row 11 ("en", 1, 41) is now row 10;
# bump its C value from 41 to 43.
$ws.Range("C10").Value = 43
